$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the D (Price) column keeps its text formatting so that
# numeric-looking strings such as "1.30" or "0.998" are not
# auto-converted into Excel numbers (which would drop trailing zeros
# or change the stored type).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '98.034.88'
$ws.Range("E2").Value = '  +3.34%  '
$ws.Range("D3").Value = '3.315.11'
$ws.Range("E3").Value = '  +7.48%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '254.21'
$ws.Range("E5").Value = '  +6.98%  '
$ws.Range("D6").Value = '627.51'
$ws.Range("E6").Value = '  +3.13%  '
$ws.Range("D7").Value = '1.30'
$ws.Range("E7").Value = '  +17.32%  '
$ws.Range("D8").Value = '0.383'
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  -0.11%  '
$ws.Range("D10").Value = '3.314.26'
$ws.Range("E10").Value = '  +7.47%  '
$ws.Range("D11").Value = '0.833'
$ws.Range("E11").Value = '  +4.42%  '
$ws.Range("E12").Value = '  +1.12%  '
$ws.Range("D13").Value = '97.527.96'
$ws.Range("E13").Value = '  +3.33%  '
$ws.Range("D14").Value = '35.46'
$ws.Range("E14").Value = '  +5.44%  '
$ws.Range("D15").Value = '0.0000244'
$ws.Range("E15").Value = '  +1.25%  '
$ws.Range("D16").Value = '3.931.26'
$ws.Range("E16").Value = '  +7.38%  '
$ws.Range("D17").Value = '5.47'
$ws.Range("E17").Value = '  +2.44%  '
$ws.Range("D18").Value = '3.322.79'
$ws.Range("E18").Value = '  +7.88%  '
$ws.Range("D19").Value = '3.55'
$ws.Range("E19").Value = '  -1.12%  '
$ws.Range("D20").Value = '14.84'
$ws.Range("E20").Value = '  +3.45%  '
$ws.Range("D21").Value = '484.19'
$ws.Range("E21").Value = '  +8.70%  '
$ws.Range("D22").Value = '6.02'
$ws.Range("E22").Value = '  +5.28%  '
$ws.Range("D23").Value = '0.0000201'
$ws.Range("E23").Value = '  +4.10%  '
$ws.Range("D24").Value = '9.21'
$ws.Range("E24").Value = '  +3.61%  '
$ws.Range("D25").Value = '5.77'
$ws.Range("E25").Value = '  +4.38%  '
$ws.Range("E26").Value = '  +4.19%  '
$ws.Range("D27").Value = '11.86'
$ws.Range("E27").Value = '  +1.09%  '
$ws.Range("D28").Value = '3.526.86'
$ws.Range("E28").Value = '  +8.68%  '
$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D29").Value = '0.278'
$ws.Range("E29").Value = '  +14.60%  '
$ws.Range("B30").Value = 'Dai'
$ws.Range("C30").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("D31").Value = '0.189'
$ws.Range("E31").Value = '  +8.14%  '
$ws.Range("D32").Value = '0.123'
$ws.Range("E32").Value = '  -9.23%  '
$ws.Range("D33").Value = '0.998'
$ws.Range("E33").Value = '  +0.13%  '
$ws.Range("D34").Value = '9.17'
$ws.Range("E34").Value = '  +2.30%  '
$ws.Range("D35").Value = '27.01'
$ws.Range("E35").Value = '  +5.52%  '
$ws.Range("D36").Value = '0.152'
$ws.Range("E36").Value = '  +0.64%  '
$ws.Range("D37").Value = '510.78'
$ws.Range("E37").Value = '  +4.31%  '
$ws.Range("D38").Value = '7.25'
$ws.Range("E38").Value = '  -2.45%  '
$ws.Range("E39").Value = '  +4.67%  '
$ws.Range("D40").Value = '24.79'
$ws.Range("E40").Value = '  +2.98%  '
$ws.Range("D41").Value = '0.448'
$ws.Range("E41").Value = '  +3.45%  '
$ws.Range("D42").Value = '3.79'
$ws.Range("E42").Value = '  -1.25%  '
$ws.Range("D43").Value = '1.25'
$ws.Range("E43").Value = '  +1.26%  '
$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value = '3.16'
$ws.Range("E45").Value = '  +0.16%  '
$ws.Range("D46").Value = '0.760'
$ws.Range("E46").Value = '  +11.88%  '
$ws.Range("D47").Value = '160.39'
$ws.Range("E47").Value = '  -0.65%  '
$ws.Range("D48").Value = '1.94'
$ws.Range("E48").Value = '  +6.16%  '
$ws.Range("D49").Value = '45.43'
$ws.Range("E49").Value = '  +3.98%  '
$ws.Range("B50").Value = 'Filecoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D50").Value = '4.51'
$ws.Range("E50").Value = '  +5.69%  '
$ws.Range("B51").Value = 'ImmutableX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D51").Value = '1.35'
$ws.Range("E51").Value = '  +5.93%  '

# Restore the default cell style on the price column so no stray
# number-format styles are left behind in the workbook.
$ws.Range("D2:D51").Style = "Normal"
